# Generate Report for Handoff
# Update the localization-status report so that the rows which just became
# "Ready for handoff" (rows 7, 9, 10, 12, 13, 14 on the language sheets) get:
#   - Priority ("E" column) set to "ht"
#   - Latest Handoff Datetime ("H" column) refreshed to the new generation time
# and the Overview sheet's "Latest HO Xliff Generate Date" ("G" column) for
# those same source files is refreshed as well.

$wb = $excel.ActiveWorkbook

$rows = @(7, 9, 10, 12, 13, 14)

# Overview sheet: column G = "Latest HO Xliff Generate Date"
$overview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $overview.Cells.Item($r, 7).Value = "2016-08-28 06:25:24"
}

# zh-cn sheet: column E = "Priority", column H = "Latest Handoff Datetime"
$zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $zhcn.Cells.Item($r, 5).Value = "ht"
    $zhcn.Cells.Item($r, 8).Value = "2016-08-28 06:25:19"
}

# de-de sheet: column E = "Priority", column H = "Latest Handoff Datetime"
$dede = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $dede.Cells.Item($r, 5).Value = "ht"
    $dede.Cells.Item($r, 8).Value = "2016-08-28 06:25:24"
}
